# Apply the "add mapping to multiple metric frameworks" edit to the
# data_dictionary workbook.
#
# In the line_items.tsv block of Sheet1:
#   - rows 41-43 (indicator / score / attribute) are removed entirely
#   - row 40 (metric_id) is renamed to metric_ids, with an updated
#     definition and a new note about disambiguating line items that map
#     to more than one metric in a given framework

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Remove the now-redundant indicator / score / attribute rows that used to
# follow metric_id in the line_items.tsv section (rows 41, 42, 43). Deleting
# them shifts every subsequent row up by three.
$ws.Rows.Item(41).Resize(3).EntireRow.Delete()

# Update the metric_id row to reflect that a line item can now map to more
# than one metric id per framework.
$ws.Range("B40").Value = "metric_ids"
$ws.Range("C40").Value = "A unique ID (or unique IDs) associated with the specified metric(s), which can be used to join to the metrics.tsv table, metrics are listed in a comma serparated list"
$ws.Range("D40").Value = "Where a single line-item corresponding to more than one metric from a given framework (e.g. JEE or SPAR), efforts where made to disambiguate to the extent possible in order to assign each line item to a maximum of one specific metric per framework. For example, a given cost will not be mapped to more than one metric of the JEE. Instead, each line item was mapped to the single most relevant metric per framework."

# Refresh the visible scroll position/selection to match the author's saved
# view (they were reviewing the edited rows around 37-41 when they saved).
$ws.Application.Goto($ws.Range("A37"), $true)
$ws.Range("C41").Select()
